$d = $word.ActiveDocument

# The document contains a spurious, fully-empty paragraph (just a run with
# color formatting and no text) sitting right after the
# "... Services registration :" paragraph and before "End of demonstration.".
# This reproduces the "empty AQL expression generates an empty line" bug;
# the fix removes that blank paragraph entirely (fixed #418).

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text.StartsWith("End of demonstration.")) {
            $target = $p
            break
        }
    }
}

if ($target -ne $null) {
    # Remove the whole paragraph, including its trailing paragraph mark,
    # so the surrounding paragraphs merge back together with no blank line
    # left behind.
    $target.Range.Delete()
} else {
    # Fallback: the empty paragraph is known to be the 2nd paragraph in
    # this document.
    $d.Paragraphs(2).Range.Delete()
}
